$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-10-30"
$ws.Range("A6").NumberFormat = "General"
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = "NAMI"
$ws.Range("D6").Value = "Cleaning"
$ws.Range("E6").Value = 500
$ws.Range("F6").Value = "TEST2"
